$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("N1").Value = "REMINDER_ROW_ID"
$ws.Range("O1").Value = "REMINDER_SNOOZE_UNTIL"
$ws.Range("P1").Value = "REMINDER_DISMISSED"

$ws.Range("N1:P1").Style = $ws.Range("M1").Style

# Row data: RowNumber -> UUID
$rowIds = @{
    2 = "11b46c4e-084f-42b1-acd9-7a800fb62cf4"
    3 = "d5caa7d8-6422-4365-ba4c-6033a58cf91a"
    4 = "4b013cf0-75ea-4221-bc28-0d683724c08a"
    5 = "ba135f90-2877-472c-b0af-fffbab774325"
    6 = "77bb20c2-2f4a-43fe-9275-a10bde9b17f9"
}

foreach ($r in 2..6) {
    $ws.Cells.Item($r, 14).Value = $rowIds[$r]
    $ws.Cells.Item($r, 15).Value = ""
    $ws.Cells.Item($r, 16).Value = $false
}
